$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44357, 50, 6000, 6500, 6200, 620),
    @(3, 44204, 80, 7000, 7500, 7188, 719),
    @(4, 44349, 60, 6000, 6500, 6250, 625),
    @(5, 44327, 60, 7000, 7500, 7250, 725),
    @(6, 44188, 80, 8000, 8500, 8250, 825),
    @(7, 44383, 60, 7500, 8000, 7750, 775),
    @(8, 44336, 60, 6000, 6500, 6250, 625),
    @(9, 44225, 60, 7500, 8000, 7750, 775),
    @(10, 44160, 100, 9000, 9500, 9250, 925),
    @(11, 44425, 60, 6500, 7000, 6750, 675),
    @(12, 44166, 100, 8000, 9000, 8500, 850),
    @(13, 44405, 80, 7500, 8000, 7688, 769),
    @(14, 44308, 100, 5000, 5500, 5250, 525),
    @(15, 44231, 70, 7500, 8000, 7714, 771),
    @(16, 44355, 50, 6000, 6500, 6300, 630),
    @(17, 44246, 60, 9000, 10000, 9500, 950),
    @(18, 44230, 60, 9000, 10000, 9500, 950),
    @(19, 44328, 60, 7000, 7500, 7250, 725),
    @(20, 44320, 50, 7000, 7500, 7200, 720),
    @(21, 44209, 80, 7500, 8000, 7688, 769),
    @(22, 44334, 60, 6500, 7000, 6750, 675),
    @(23, 44299, 100, 8000, 9000, 8500, 850),
    @(24, 44292, 50, 10000, 11000, 10600, 1060),
    @(25, 44362, 50, 6000, 6500, 6300, 630),
    @(26, 44365, 50, 6000, 6500, 6200, 620),
    @(27, 44433, 100, 7000, 7500, 7250, 725),
    @(28, 44421, 100, 7000, 7500, 7250, 725),
    @(29, 44194, 100, 8000, 9000, 8500, 850),
    @(30, 44264, 50, 8000, 8500, 8200, 820),
    @(31, 44316, 100, 6000, 6500, 6250, 625),
    @(32, 44351, 50, 6000, 6500, 6300, 630),
    @(33, 44273, 80, 7000, 8000, 7500, 750),
    @(34, 44342, 50, 6000, 6500, 6300, 630),
    @(35, 44313, 60, 6000, 6500, 6250, 625),
    @(36, 44428, 50, 7500, 8000, 7800, 780),
    @(37, 44435, 100, 7000, 7500, 7250, 725),
    @(38, 44376, 100, 6000, 6500, 6250, 625),
    @(39, 44279, 60, 7500, 8000, 7750, 775),
    @(40, 44238, 100, 8000, 8500, 8250, 825),
    @(41, 44399, 60, 9000, 10000, 9500, 950)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P: Precio $/Kg
}
